$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - names rotated
$ws.Range("A1").Value = "Søren"
$ws.Range("C1").Value = "Mads"
$ws.Range("E1").Value = "Kim"
$ws.Range("G1").Value = "Emil"

# Column A (Søren's picks)
$ws.Range("A3").Value = "Tottenham"
$ws.Range("A4").Value = "Freiburg"
$ws.Range("A5").Value = "Atlético Madrid"
$ws.Range("A6").Value = "Inter"
$ws.Range("A7").Value = "Brøndby"

# Column C (Mads' picks)
$ws.Range("C2").Value = "Liverpool"
$ws.Range("C3").Value = "Bayern"
$ws.Range("C4").Value = "Union Berlin"
$ws.Range("C5").Value = "Sevilla"
$ws.Range("C6").Value = "Napoli"
$ws.Range("C7").Value = "FC Midtjylland"

# Column E (Kim's picks)
$ws.Range("E3").Value = "Dortmund"
$ws.Range("E4").Value = "1. FC Köln"
$ws.Range("E5").Value = "Barcelona"
$ws.Range("E6").Value = "Juventus"

# Column G (Emil's picks)
$ws.Range("G2").Value = "Man City"
$ws.Range("G4").Value = "RB Leipzig"
$ws.Range("G5").Value = "Real Madrid"
$ws.Range("G6").Value = "Milan"
$ws.Range("G7").Value = "AaB"
